$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list (price/volume columns) with refreshed quote data.
# Leading "'" forces cells whose new text is numeric-looking ("561.04")
# to stay stored as text, matching the original column's text formatting.
$ws.Range("D2").Value = "64.375.82"
$ws.Range("D3").Value = "3.088.06"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'561.04"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").Value = "'145.20"
$ws.Range("E6").Value = "  +3.74%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.084.18"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").Value = "'6.14"
$ws.Range("E11").Value = "  -4.36%  "
$ws.Range("E12").Value = "  +4.25%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").Value = "'35.24"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "3.582.84"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "64.382.76"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "3.086.91"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").Value = "'480.70"
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("D21").Value = "'13.93"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").Value = "'0.674"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  +4.71%  "
$ws.Range("D24").Value = "'13.87"
$ws.Range("E24").Value = "  +10.04%  "
$ws.Range("D25").Value = "'81.23"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D28").Value = "'8.03"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("E29").Value = "  +4.60%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'26.27"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").Value = "'1.14"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").Value = "'5.61"
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("D35").Value = "'6.22"
$ws.Range("E35").Value = "  +4.09%  "
$ws.Range("D36").Value = "'55.82"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").Value = "'3.04"
$ws.Range("E37").Value = "  +17.44%  "
$ws.Range("D38").Value = "'458.80"
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "2.981.08"
$ws.Range("E41").Value = "  -2.87%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("D44").Value = "'28.05"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("D46").Value = "'2.15"
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("D49").Value = "'121.06"
$ws.Range("E49").Value = "  +3.50%  "
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("E51").Value = "  +0.72%  "
